$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.473.26"
$ws.Range("E2").Value = "  -0.62%  "

# Row 3
$ws.Range("D3").Value = "2.625.64"
$ws.Range("E3").Value = "  +0.44%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.99%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.554"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.87%  "

# Row 9
$ws.Range("D9").Value = "2.624.26"
$ws.Range("E9").Value = "  +0.42%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.122"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.51%  "

# Row 11
$ws.Range("E11").Value = "  +0.31%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.00%  "

# Row 13
$ws.Range("E13").Value = "  -1.87%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.48%  "

# Row 15
$ws.Range("D15").Value = "3.088.30"
$ws.Range("E15").Value = "  -0.06%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.06%  "

# Row 17
$ws.Range("D17").Value = "67.352.12"
$ws.Range("E17").Value = "  -0.47%  "

# Row 18
$ws.Range("D18").Value = "2.629.82"
$ws.Range("E18").Value = "  +0.58%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "370.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.80%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.72%  "

# Row 21
$ws.Range("E21").Value = "  -3.55%  "

# Row 22
$ws.Range("E22").Value = "  -0.41%  "

# Row 23
$ws.Range("B23").Value = "NEARProtocol"
$ws.Range("C23").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.58%  "

# Row 24
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.87%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "

# Row 26
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.93%  "

# Row 27
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "66.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.59%  "

# Row 28
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.740.37"
$ws.Range("E28").Value = "  -0.06%  "

# Row 29
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "587.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.16%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.22%  "

# Row 31
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0000101"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.57%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.34%  "

# Row 33
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.11%  "

# Row 34
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.87%  "

# Row 35
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.123"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.21%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.81%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.44%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.64%  "

# Row 40
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.96%  "

# Row 41
$ws.Range("E41").Value = "  -0.62%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.70%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.71%  "

# Row 44
$ws.Range("E44").Value = "  +2.24%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "155.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.12%  "

# Row 47
$ws.Range("D47").Value = "0.0₆0298"
$ws.Range("E47").Value = "  +0.53%  "

# Row 48
$ws.Range("E48").Value = "  -0.33%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.45%  "

# Row 50
$ws.Range("E50").Value = "  -1.28%  "

# Row 51
$ws.Range("E51").Value = "  +0.73%  "
